$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the source data which stores
# prices as literal strings (e.g. "27.986.23") rather than numeric values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.035.85"
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("D3").Value = "1.726.29"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "219.42"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").Value = "0.525"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "24.18"
$ws.Range("E8").Value = "  +13.78%  "
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").Value = "0.0903"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").Value = "1.970.77"
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("D13").Value = "1.723.42"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("E15").Value = "  +6.02%  "
$ws.Range("D16").Value = "67.90"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").Value = "27.966.47"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").Value = "243.44"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").Value = "0.0₃0759"
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").Value = "7.92"
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  +4.45%  "
$ws.Range("D23").Value = "9.80"
$ws.Range("E23").Value = "  +4.76%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "149.21"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").Value = "7.54"
$ws.Range("E26").Value = "  +4.28%  "
$ws.Range("D27").Value = "16.83"
$ws.Range("E27").Value = "  +2.83%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  +2.88%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").Value = "3.46"
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("D33").Value = "3.29"
$ws.Range("E33").Value = "  +3.16%  "
$ws.Range("D34").Value = "1.487.57"
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("D36").Value = "0.963"
$ws.Range("E36").Value = "  +4.43%  "
$ws.Range("D37").Value = "0.612"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").Value = "0.0176"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("E41").Value = "  +5.59%  "
$ws.Range("E42").Value = "  +4.14%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").Value = "1.874.75"
$ws.Range("E45").Value = "  +3.02%  "
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").Value = "1.77"
$ws.Range("E47").Value = "  +13.10%  "
$ws.Range("D48").Value = "91.86"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "0.0₆0110"
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.106"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "8.23"
$ws.Range("E51").Value = "  +2.17%  "

# Restore default (no explicit) style on column D now that the text values are set,
# so the cells match the original unstyled text-cell layout.
$ws.Range("D2:D51").Style = "Normal"
